# Updated Efna1 -> Epha1 NATMI edge table (new run with M2 cluster, averaged over replicates)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,("ECs", "Efna1", "Epha1", "ECs", 3, 1, 23.82411166666667, 71.472335, 0.8811513593020099, 0.8811513593020099, 3, 1, 3.746298666666667, 11.238896, 0.2944715732214049, 0.294471573221405, 89.25223777135112, 803.2701399421601, 0.2594740270198423, 0.2594740270198423)
    ,("ECs", "Efna1", "Epha1", "FAPs", 3, 1, 23.82411166666667, 71.472335, 0.8811513593020099, 0.8811513593020099, 3, 1, 4.370261666666667, 13.110785, 0.3435171466234404, 0.3435171466234404, 104.1176019592195, 937.058417632975, 0.3026906006907923, 0.3026906006907923)
    ,("ECs", "Efna1", "Epha1", "M2", 3, 1, 23.82411166666667, 71.472335, 0.8811513593020099, 0.8811513593020099, 3, 1, 1.837969333333333, 5.513908, 0.1444705212467569, 0.1444705212467569, 43.78798663724222, 394.09187973518, 0.1273003961756497, 0.1273003961756497)
    ,("ECs", "Efna1", "Epha1", "sCs", 3, 1, 23.82411166666667, 71.472335, 0.8811513593020099, 0.8811513593020099, 3, 1, 2.767576666666667, 8.30273, 0.2175407589083977, 0.2175407589083977, 65.93505555272777, 593.41549997455, 0.1916863354157255, 0.1916863354157255)
    ,("FAPs", "Efna1", "Epha1", "ECs", 3, 1, 2.471212, 7.413636, 0.09139949658522162, 0.09139949658522162, 3, 1, 3.746298666666667, 11.238896, 0.2944715732214049, 0.294471573221405, 9.257898220650667, 83.32108398585601, 0.02691455355109464, 0.02691455355109464)
    ,("FAPs", "Efna1", "Epha1", "FAPs", 3, 1, 2.471212, 7.413636, 0.09139949658522162, 0.09139949658522162, 3, 1, 4.370261666666667, 13.110785, 0.3435171466234404, 0.3435171466234404, 10.79984307380667, 97.19858766426, 0.03139729426977422, 0.03139729426977422)
    ,("FAPs", "Efna1", "Epha1", "M2", 3, 1, 2.471212, 7.413636, 0.09139949658522162, 0.09139949658522162, 3, 1, 1.837969333333333, 5.513908, 0.1444705212467569, 0.1444705212467569, 4.542011872165333, 40.878106849488, 0.01320453291335814, 0.01320453291335814)
    ,("FAPs", "Efna1", "Epha1", "sCs", 3, 1, 2.471212, 7.413636, 0.09139949658522162, 0.09139949658522162, 3, 1, 2.767576666666667, 8.30273, 0.2175407589083977, 0.2175407589083977, 6.839268669586667, 61.55341802628001, 0.01988311585099462, 0.01988311585099462)
    ,("M2", "Efna1", "Epha1", "ECs", 1, 0.3333333333333333, 0.02836866666666667, 0.085106, 0.001049234890461559, 0.001049234890461559, 3, 1, 3.746298666666667, 11.238896, 0.2944715732214049, 0.294471573221405, 0.1062774981084444, 0.956497482976, 0.0003089698488730038, 0.0003089698488730039)
    ,("M2", "Efna1", "Epha1", "FAPs", 1, 0.3333333333333333, 0.02836866666666667, 0.085106, 0.001049234890461559, 0.001049234890461559, 3, 1, 4.370261666666667, 13.110785, 0.3435171466234404, 0.3435171466234404, 0.1239784964677778, 1.11580646821, 0.0003604301757091129, 0.0003604301757091129)
    ,("M2", "Efna1", "Epha1", "M2", 1, 0.3333333333333333, 0.02836866666666667, 0.085106, 0.001049234890461559, 0.001049234890461559, 3, 1, 1.837969333333333, 5.513908, 0.1444705212467569, 0.1444705212467569, 0.05214073936088889, 0.469266654248, 0.0001515835115352653, 0.0001515835115352653)
    ,("M2", "Efna1", "Epha1", "sCs", 1, 0.3333333333333333, 0.02836866666666667, 0.085106, 0.001049234890461559, 0.001049234890461559, 3, 1, 2.767576666666667, 8.30273, 0.2175407589083977, 0.2175407589083977, 0.07851245993111111, 0.70661213938, 0.0002282513543441771, 0.0002282513543441771)
    ,("sCs", "Efna1", "Epha1", "ECs", 3, 1, 0.7137869999999999, 2.141361, 0.02639990922230694, 0.02639990922230694, 3, 1, 3.746298666666667, 11.238896, 0.2944715732214049, 0.294471573221405, 2.674059286384, 24.066533577456, 0.007774022801595002, 0.007774022801595002)
    ,("sCs", "Efna1", "Epha1", "FAPs", 3, 1, 0.7137869999999999, 2.141361, 0.02639990922230694, 0.02639990922230694, 3, 1, 4.370261666666667, 13.110785, 0.3435171466234404, 0.3435171466234404, 3.119435964265, 28.074923678385, 0.00906882148716473, 0.00906882148716473)
    ,("sCs", "Efna1", "Epha1", "M2", 3, 1, 0.7137869999999999, 2.141361, 0.02639990922230694, 0.02639990922230694, 3, 1, 1.837969333333333, 5.513908, 0.1444705212467569, 0.1444705212467569, 1.311918616532, 11.807267548788, 0.003814008646213747, 0.003814008646213747)
    ,("sCs", "Efna1", "Epha1", "sCs", 3, 1, 0.7137869999999999, 2.141361, 0.02639990922230694, 0.02639990922230694, 3, 1, 2.767576666666667, 8.30273, 0.2175407589083977, 0.2175407589083977, 1.97546024617, 17.77914221553, 0.00574305628733346, 0.00574305628733346)
)

$numRows = $rows.Count
$numCols = $rows[0].Count
$data = New-Object 'object[,]' $numRows, $numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $data[$r, $c] = $rows[$r][$c]
    }
}

$ws.Range("A2:T17").Value = $data
